$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (dates as text to match existing column A formatting,
# numeric rates in columns B-E).
$rows = @(
    @(87, "2025-09-05", 0.04217290000000001, 0.0412319,            0.0396676, 0.0371438),
    @(88, "2025-09-08", 0.0419217,            0.0406989,            0.0388309, 0.0360736),
    @(89, "2025-09-09", 0.0417535,            0.04055950000000001, 0.0386492, 0.0359297),
    @(90, "2025-09-10", 0.041739,             0.0405836,            0.038789,  0.0361713),
    @(91, "2025-09-11", 0.0415016,            0.04037810000000001, 0.0386433, 0.0361192)
)

foreach ($row in $rows) {
    $r = $row[0]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[1]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
